$wb = $excel.ActiveWorkbook

# Rename "Step 3" to "Step 4"
$ws4 = $wb.Worksheets.Item("Step 3")
$ws4.Name = "Step 4"

# Set B11 on "Setp 2" to "Yes"
$wsSetp2 = $wb.Worksheets.Item("Setp 2")
$wsSetp2.Range("B11").Value = "Yes"
